$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.220.85"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.809.56"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.74"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.16"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0720"
$ws.Range("E10").Value = "  +8.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "2.075.98"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "1.818.77"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.91"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.637"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "34.201.50"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.29"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.59"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.71"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  +6.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.96"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.24"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.16"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.67"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0532"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.58"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "1.433.38"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.637"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.956"
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.99"
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.97"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0498"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.05"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.967.56"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.15"
$ws.Range("E48").Value = "  +6.92%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.996"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.87"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("E51").Value = "  +6.76%  "
